$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Flip the direction input cell from "ltr" to "rtl".
$ws.Range("C3").Value = "rtl"

# The "prev" column (D) holds hand-maintained axle locations that cascade
# off of the spacing column (C) and the direction flag above. These aren't
# formulas, so they need to be re-keyed by hand to the values that match
# the new ("rtl") direction -- mirroring D(r) = D(r-1) +/- C(r).
$ws.Range("D6").Value = 83
$ws.Range("D7").Value = 88
$ws.Range("D8").Value = 93
$ws.Range("D9").Value = 98
$ws.Range("D10").Value = 107
$ws.Range("D11").Value = 112
$ws.Range("D12").Value = 118
$ws.Range("D13").Value = 123
$ws.Range("D14").Value = 131
$ws.Range("D15").Value = 139
$ws.Range("D16").Value = 144
$ws.Range("D17").Value = 149
$ws.Range("D18").Value = 154
$ws.Range("D19").Value = 163
$ws.Range("D20").Value = 168
$ws.Range("D21").Value = 174
$ws.Range("D22").Value = 179

# Restore the view state left behind by the edit: scrolled so column AA is
# at the top-left, with AO32 as the active selection.
$excel.ActiveWindow.ScrollColumn = 27
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("AO32").Select()
